$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row
# D-column values are forced to Text format (then style restored) so that
# numeric-looking strings (e.g. "1.000", "19.00", "0.06539") keep their exact
# textual representation instead of being auto-converted to numbers by Excel.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.203.65"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -1.49%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.862.75"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -2.17%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9992"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  -0.36%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "234.37"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -2.27%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -0.29%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4675"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -1.31%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2834"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -1.05%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06539"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -1.90%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.77"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +5.78%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07862"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +0.75%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "96.56"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -4.37%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.861.90"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -2.25%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.123"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -1.26%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6721"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -1.15%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "279.54"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -2.22%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.205.47"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -1.53%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9979"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -0.41%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.489"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.19%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.66"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.65%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.104.78"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -3.19%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.000007257"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  -3.25%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -0.26%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.153"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -1.85%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.309"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -0.62%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "164.92"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -1.38%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.00"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -1.99%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.913"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -6.13%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.351"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -2.22%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09573"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -3.91%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.416"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -2.58%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.467"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -3.26%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.111"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -3.33%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.04698"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.26%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7022"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -2.95%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.098"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("E37").Value = "  -0.39%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.01874"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -2.05%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.348"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -7.15%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.521"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -3.03%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "73.14"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -1.49%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.937"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -2.85%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.8477"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -2.18%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4180"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -2.39%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -0.32%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "103.83"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -1.37%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "990.82"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -1.49%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.161"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -2.87%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.260"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.31%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "34.02"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -1.44%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1135"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -4.13%  "
